$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the entire B2:D9 block to 0
$ws.Range("B2:D9").Value = 0

# Apply the two specific non-zero overrides from the diff
$ws.Range("C4").Value = 0.7840954128749528
$ws.Range("C8").Value = -0.6615990660246527
